# ANOVA results: fill in significance codes for AG, AP, BG, BX, CBH, LAP
# (NAG and PPO are left for later, per commit message)

$wb = $excel.ActiveWorkbook
$wsVmax = $wb.Worksheets.Item("Vmax")
$wsKm   = $wb.Worksheets.Item("Km")

# --- Vmax sheet ---
# row 4 = BG
$wsVmax.Range("B4").Value = "*"
$wsVmax.Range("C4").Value = "***"
$wsVmax.Range("D4").Value = "o"
$wsVmax.Range("E4").Value = "o"
$wsVmax.Range("F4").Value = "o"
$wsVmax.Range("G4").Value = "o"
$wsVmax.Range("H4").Value = "o"
# row 5 = BX
$wsVmax.Range("B5").Value = "o"
$wsVmax.Range("C5").Value = "o"
$wsVmax.Range("D5").Value = "o"
$wsVmax.Range("E5").Value = "o"
$wsVmax.Range("F5").Value = "o"
$wsVmax.Range("G5").Value = "o"
$wsVmax.Range("H5").Value = "o"
# row 6 = CBH
$wsVmax.Range("B6").Value = "***"
$wsVmax.Range("C6").Value = "***"
$wsVmax.Range("D6").Value = "***"
$wsVmax.Range("E6").Value = "*"
$wsVmax.Range("F6").Value = "***"
$wsVmax.Range("G6").Value = "o"
$wsVmax.Range("H6").Value = "o"
# row 7 = LAP
$wsVmax.Range("B7").Value = "o"
$wsVmax.Range("C7").Value = "o"
$wsVmax.Range("D7").Value = "o"
$wsVmax.Range("E7").Value = "o"
$wsVmax.Range("F7").Value = "o"
$wsVmax.Range("G7").Value = "o"
$wsVmax.Range("H7").Value = "o"

# --- Km sheet ---
# row 3 = AP, Precipitation column updated
$wsKm.Range("D3").Value = "*"
# row 4 = BG
$wsKm.Range("B4").Value = "o"
$wsKm.Range("C4").Value = "o"
$wsKm.Range("D4").Value = "o"
$wsKm.Range("E4").Value = "o"
$wsKm.Range("F4").Value = "o"
$wsKm.Range("G4").Value = "o"
$wsKm.Range("H4").Value = "o"
# row 5 = BX
$wsKm.Range("B5").Value = "o"
$wsKm.Range("C5").Value = "*"
$wsKm.Range("D5").Value = "o"
$wsKm.Range("E5").Value = "o"
$wsKm.Range("F5").Value = "o"
$wsKm.Range("G5").Value = "o"
$wsKm.Range("H5").Value = "o"
# row 6 = CBH
$wsKm.Range("B6").Value = "***"
$wsKm.Range("C6").Value = "***"
$wsKm.Range("D6").Value = "o"
$wsKm.Range("E6").Value = "*"
$wsKm.Range("F6").Value = "***"
$wsKm.Range("G6").Value = "o"
$wsKm.Range("H6").Value = "o"
# row 7 = LAP
$wsKm.Range("B7").Value = "o"
$wsKm.Range("C7").Value = "***"
$wsKm.Range("D7").Value = "o"
$wsKm.Range("E7").Value = "o"
$wsKm.Range("F7").Value = "o"
$wsKm.Range("G7").Value = "o"
$wsKm.Range("H7").Value = "o"

# --- Selections (match final cursor positions from the diff) ---
$wsVmax.Range("B8").Select()

$wsKm.Activate()
$wsKm.Range("B7").Select()
